# Apply the "tasks" edit described by the commit:
#   Earcons - TTS spell incoming number - Deal with private numbers
#
# This rewrites the contents of the task table on sheet "Folha1":
#  - refines a couple of task titles
#  - marks the Earcons / TTS Spell / private-number tasks as Done, owned by Hugo
#  - adds a new row for "Debug todas as combinacoes de receber chamada / em chamada"
#  - splits the old "Debug todas as combinacoes de chamadas" task into an
#    outgoing-call debug task (existing row) and the new incoming-call one

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# --- Row 2: pop-up menu task gets a clarifying suffix -----------------------
$ws.Range("A2").Value = "Lidar com pop-up menu após chamada (só vodafone?)"

# --- Row 3: unchanged (Always visible) --------------------------------------

# --- Row 4: "Debug all call combinations" becomes the outgoing-call variant,
#            now owned by Hugo -----------------------------------------------
$ws.Range("A4").Value = "Debug todas as combinações de fazer chamada "
$ws.Range("C4").Value = "Hugo"

# --- Row 5: "TTS spell function" row is replaced in place by the item that
#            used to be row 6 (Usar lista de contactos em incoming call) ----
$ws.Range("A5").Value = "Usar lista de contactos em incoming call"

# --- Row 6: "Usar lista de contactos..." row is replaced in place by
#            "Aplicar earcons", now Done and owned by Hugo, and it inherits
#            the "Earcons em acções..." description that used to sit on
#            the "Aplicar earcons" row (old row 7), including its italic
#            formatting -----------------------------------------------------
$ws.Range("A6").Value = "Aplicar earcons"
$ws.Range("C6").Value = "Hugo"
$ws.Range("D6").Value = "Done"
$ws.Range("G6").Value = "Earcons em acções de select e voltar"
$ws.Range("G6").Font.Italic = $true

# --- Row 7: "Aplicar earcons" row is replaced in place by "TTS Spell",
#            Done, owned by Hugo, description cleared (content+format) -----
$ws.Range("A7").Value = "TTS Spell"
$ws.Range("C7").Value = "Hugo"
$ws.Range("D7").Value = "Done"
$ws.Range("G7").Clear() | Out-Null

# --- Row 8: "Increase InCall TTS volume" is replaced by the new private
#            number task, Done, owned by Hugo, priority raised to 3 ---------
$ws.Range("A8").Value = "Lidar com numero privado (incoming number = null)"
$ws.Range("B8").Value = 3
$ws.Range("C8").Value = "Hugo"
$ws.Range("D8").Value = "Done"

# --- Rows 9-12: unchanged task names, but now all owned by Hugo ------------
$ws.Range("C9").Value = "Hugo"
$ws.Range("C10").Value = "Hugo"
$ws.Range("C11").Value = "Hugo"
$ws.Range("C12").Value = "Hugo"

# --- Row 13 (new): incoming-call / in-call debug task, Done, owned by Hugo -
$ws.Range("A13").Value = "Debug todas as combinações de receber chamada / em chamada"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "Hugo"
$ws.Range("D13").Value = "Done"

# Keep the autofilter / sort range in sync with the new last row and
# re-apply the existing sort (by Status, descending) so sortState covers
# the extended range A2:G13.
$ws.Range("A1:G13").Sort($ws.Range("D1"), 2, $null, $null, 1, $null, 1, 1) | Out-Null

# Match the author's final cursor position.
$ws.Range("A5").Select() | Out-Null
